# Update "想去人数" (want-to-go count) figures in the F column on both the
# "展览" (sheet1) and "全部类型" (sheet4) worksheets.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 243
    3  = 268
    6  = 274
    7  = 6634
    9  = 74
    10 = 116
    11 = 78
    12 = 37
    13 = 10
    15 = 18
    16 = 213
    17 = 556
    18 = 56
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
